$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B53: change from inline-string "2" to numeric 2
$ws.Range("B53").Value = 2

# Add new row 54
$ws.Range("A54").Value = "Ruilin"
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = "4"
$ws.Range("C54").Value = "interested to hear more"
$ws.Range("D54").Value = "DIS"
$ws.Range("E54").Value = "WRI"
$ws.Range("F54").Value = "a0a400ab-cd67-43a0-98e0-d641a379b0a8"
$ws.Range("G54").Value = "B1QRgziT-_annotated.xlsx"
$ws.Range("H54").Value = "I am also interested to hear more about the semantics of the spectral norm of this object (flattened filterbank), which Ian asked about below."
